# build: ajust prompt ai
#
# 1. Add a new row (47) to the "geral" sheet with a new MI31 transaction entry.
# 2. Move the active/selected tab from "erros" (sheet 3) to "geral" (sheet 1),
#    with the active cell on the newly added row (D47).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("geral")

# New row of data (transacao / texto / processo / descricao)
$ws.Range("A47").Value = "MI31"
$ws.Range("B47").Value = "Batch Input: criar doc.invent."
$ws.Range("C47").Value = "inventario"
$ws.Range("D47").Value = "criação de documento de inventário"

# Make "geral" the active sheet and select/scroll to the newly added cell.
$ws.Activate()
$ws.Range("D47").Select()

# Best-effort: scroll the view so row 34 is the top visible row (matches the
# original author's window position, sheetView topLeftCell="A34"). Wrapped in
# try/catch since not every host exposes window scroll state.
try {
    $excel.ActiveWindow.ScrollRow = 34
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
